# "Arm based data by timepoint" sheet: insert a new "Sample Size at
# Timepoint" column right before the existing Mean / Standard Deviation /
# Standard Error columns (which shift one column to the right).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arm based data by timepoint")

$ws.Columns.Item(7).Insert()
$ws.Cells.Item(1, 7).Value = "Sample Size at Timepoint"

# Leave the cursor where the author left it while working on this sheet ...
$ws.Range("G5").Select() | Out-Null

# ... then return focus to "Contrasts between timepoints", which is the
# tab that was active/selected when the workbook was last saved.
$ws4 = $wb.Worksheets.Item("Contrasts between timepoints")
$ws4.Activate() | Out-Null
$ws4.Range("D1").Select() | Out-Null
